$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the steel ("S") description in the industrial mapping text by
# removing the "RME/" fragment from the two "S/LFM+..." lines.
$cell = $ws.Range("B2")
$cell.Replace("RME/", "") | Out-Null

# Wrap the multi-line text and let the row grow to hold it.
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.6

# Restore the on-screen selection over the updated data.
$ws.Range("B2:B12").Select() | Out-Null
